$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 27; A = 42602.513599537036; C = 2793; D = 151; E = 8; F = 2; G = 1;  H = 66;  I = 33; J = 3; K = 3; L = 49; M = 49 },
    @{ Row = 28; A = 42602.516712962963; C = 3127; D = 120; E = 5; F = 2; G = 0;  H = 100; I = 0;  J = 3; K = 2; L = 60; M = 40 },
    @{ Row = 29; A = 42602.524062500001; C = 3092; D = 120; E = 5; F = 2; G = 0;  H = 100; I = 0;  J = 3; K = 2; L = 60; M = 40 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "Named"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}
